# "added test for Line"
# A new timesheet entry (2014-04-16, 09:30 -> 09:45) is inserted right after the
# last data row (row 128), pushing the blank separator row and the three
# summary rows ("sum [min]", "sum [h]", "sum [working weeks]") down by one row.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Duplicate the last data row's structure (styles/column layout) into a new
# row 129, shifting the old row 129 (blank separator) and the summary rows
# below it down to rows 130-133.
$ws.Range("A128:G128").Copy()
$null = $ws.Range("A129:G129").Insert(-4162)  # xlShiftDown

# Fill in the new timesheet entry.
$ws.Range("A129").Value = 2014
$ws.Range("B129").Value = 4
$ws.Range("C129").Value = 16
$ws.Range("D129").Value = 0.39583333333333331   # 09:30
$ws.Range("E129").Value = 0.40625                # 09:45
$ws.Range("F129").Formula = "=(E129-D129)*24*60"
$ws.Range("G129").Formula = "=F129/60"

# Recalculate everything (sum/avg rows below pick up the new row automatically
# since their SUM(...) ranges grow along with the inserted row).
$excel.CalculateFull()

# Match the author's final selection.
$null = $ws.Range("H131").Select()
